$d = $word.ActiveDocument

$replacements = @(
    @{old = "91×92="; new = "68×64="},
    @{old = "84×16="; new = "47×13="},
    @{old = "13×79="; new = "41×33="},
    @{old = "34×90="; new = "53×56="},
    @{old = "59×31="; new = "43×99="},
    @{old = "32×68="; new = "76×99="},
    @{old = "23×37="; new = "58×82="},
    @{old = "42×85="; new = "47×56="},
    @{old = "46×40="; new = "93×51="},
    @{old = "81×41="; new = "31×61="},
    @{old = "43×73="; new = "37×57="},
    @{old = "86×67="; new = "80×48="},
    @{old = "90×11="; new = "32×14="},
    @{old = "36×35="; new = "38×89="},
    @{old = "80×65="; new = "74×32="},
    @{old = "97×44="; new = "95×39="},
    @{old = "37×56="; new = "36×78="},
    @{old = "42×86="; new = "66×44="},
    @{old = "88×66="; new = "36×59="},
    @{old = "45×64="; new = "17×76="},
    @{old = "70×58="; new = "83×20="},
    @{old = "84×12="; new = "48×60="},
    @{old = "65×95="; new = "54×37="},
    @{old = "97×99="; new = "24×53="},
    @{old = "86×96="; new = "54×73="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $false, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
